# Edit slide 33 ("Scanning a Hexadecimal Literal") of the lexical-analysis
# deck so that the hexadecimal-literal scanning example is turned into a
# binary-literal scanning example.
#
# TextRange.Find(FindWhat, After, MatchCase, WholeWords) -- WholeWords must
# be $false here since several of the target strings start/end mid-word or
# with punctuation/spaces (e.g. " contains \"0X\"").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(33)

# --- Title placeholder -----------------------------------------------------
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$found = $titleRange.Find("Scanning a Hexadecimal Literal", 0, $true, $false)
$found.Text = "Scanning a Binary Literal"

# --- Content placeholder ----------------------------------------------------
$content = $s.Shapes.Item(2)
$tr = $content.TextFrame.TextRange

# Each Find/replace below targets the *entire* text of a single run (not
# just the changed substring) so the run isn't split into extra pieces and
# the existing a:rPr formatting is preserved unchanged on one run.

# scanHexLiteral -> scanBinaryLiteral  (function name in the "private fun" line)
$found = $tr.Find("scanHexLiteral", 0, $true, $false)
$found.Text = "scanBinaryLiteral"

# " contains ""0X""" -> " contains ""0B"""
$found = $tr.Find(" contains ""0X""", 0, $true, $false)
$found.Text = " contains ""0B"""

# [1] == 'X') -> [1] == 'B')
$found = $tr.Find("[1] == 'X')", 0, $true, $false)
$found.Text = "[1] == 'B')"

# "hex digit" comment -> "binary digit" comment
$found = $tr.Find("    // check that the next character is a hex digit", 0, $true, $false)
$found.Text = "    // check that the next character is a binary digit"

# CharUtil.isHexDigit -> CharUtil.isBinaryDigit (two occurrences)
$found = $tr.Find("CharUtil.isHexDigit", 0, $true, $false)
$found.Text = "CharUtil.isBinaryDigit"

$found = $tr.Find("CharUtil.isHexDigit", $found.Start + $found.Length, $true, $false)
$found.Text = "CharUtil.isBinaryDigit"

# "Improperly formed hexadecimal literal." -> "Improperly formed binary literal."
$found = $tr.Find("        throw error(""Improperly formed hexadecimal literal."")", 0, $true, $false)
$found.Text = "        throw error(""Improperly formed binary literal."")"
